$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{Row=2; D='67.039.17'; E='  +0.40%  '}
    @{Row=3; D='3.821.74'; E='  -0.75%  '}
    @{Row=4; E='  -0.24%  '}
    @{Row=5; D='448.40'; E='  +6.17%  '}
    @{Row=6; D='146.51'; E='  +13.61%  '}
    @{Row=7; E='  +3.31%  '}
    @{Row=8; E='  +0.02%  '}
    @{Row=9; E='  +2.93%  '}
    @{Row=10; D='0.154'; E='  -5.51%  '}
    @{Row=11; D='0.0000316'; E='  -9.69%  '}
    @{Row=12; D='43.90'; E='  +10.05%  '}
    @{Row=13; D='10.40'; E='  +4.38%  '}
    @{Row=14; D='4.396.95'; E='  -1.47%  '}
    @{Row=15; D='14.75'; E='  -8.74%  '}
    @{Row=16; E='  -0.20%  '}
    @{Row=17; D='3.809.28'; E='  -0.82%  '}
    @{Row=18; D='20.01'; E='  +2.89%  '}
    @{Row=19; E='  +7.29%  '}
    @{Row=20; D='67.089.46'; E='  +0.20%  '}
    @{Row=21; D='420.98'; E='  +4.40%  '}
    @{Row=22; D='14.67'; E='  +3.34%  '}
    @{Row=23; D='3.27'; E='  +9.93%  '}
    @{Row=24; D='86.18'; E='  +2.80%  '}
    @{Row=25; B='PancakeSwap'; C='https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'; D='3.47'; E='  +9.19%  '}
    @{Row=26; B='EthereumClassic'; C='https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'; D='37.48'; E='  +1.85%  '}
    @{Row=27; D='5.55'; E='  -5.20%  '}
    @{Row=28; D='9.80'; E='  +4.60%  '}
    @{Row=29; D='9.39'; E='  +27.60%  '}
    @{Row=30; D='733.86'; E='  +2.92%  '}
    @{Row=31; D='13.77'; E='  +12.22%  '}
    @{Row=32; E='  +10.42%  '}
    @{Row=33; E='  -0.19%  '}
    @{Row=34; D='44.50'; E='  +18.69%  '}
    @{Row=35; E='  +7.23%  '}
    @{Row=36; D='56.65'; E='  +3.59%  '}
    @{Row=37; B='Dai'; C='https://coinranking.com/coin/MoTuySvg7+dai-dai'; D='1.00'; E='  +0.07%  '}
    @{Row=38; B='NEARProtocol'; C='https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'; D='5.52'; E='  +24.67%  '}
    @{Row=39; E='  +6.73%  '}
    @{Row=40; D='2.89'; E='  +0.06%  '}
    @{Row=41; D='0.339'; E='  +16.92%  '}
    @{Row=42; D='0.0₃0676'; E='  -11.47%  '}
    @{Row=43; B='Stellar'; C='https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'; D='0.141'; E='  +4.93%  '}
    @{Row=44; B='FirstDigitalUSD'; C='https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'; D='1.00'; E='  -0.17%  '}
    @{Row=45; B='Fetch.AI'; C='https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'; D='2.56'; E='  +25.48%  '}
    @{Row=46; D='3.24'; E='  +2.46%  '}
    @{Row=47; D='3.40'; E='  +3.20%  '}
    @{Row=48; D='2.13'; E='  +4.47%  '}
    @{Row=49; D='144.99'; E='  +1.01%  '}
    @{Row=50; E='  +5.09%  '}
    @{Row=51; E='  +4.43%  '}
)

foreach ($r in $rows) {
    $row = $r.Row

    if ($r.ContainsKey('B')) {
        $ws.Cells.Item($row, 2).Value = $r.B
    }
    if ($r.ContainsKey('C')) {
        $ws.Cells.Item($row, 3).Value = $r.C
    }
    if ($r.ContainsKey('D')) {
        $dcell = $ws.Cells.Item($row, 4)
        if ($r.D -match '^[0-9]+(\.[0-9]+)?$') {
            # Force as text to avoid Excel auto-converting numeric-looking
            # strings (e.g. "448.40") into floating point numbers.
            $dcell.NumberFormat = "@"
            $dcell.Value = $r.D
            $dcell.Style = "Normal"
        } else {
            $dcell.Value = $r.D
        }
    }
    if ($r.ContainsKey('E')) {
        $ws.Cells.Item($row, 5).Value = $r.E
    }
}
